# Week 9 tutorial attempt, Q1 half way done.
# Rewrites column B (rows 2-100) with the updated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3,2,1,2,1,1,4,1,4,2,2,3,5,5,1,2,3,3,5,1,3,2,5,1,2,5,2,1,5,1,1,2,1,2,2,3,4,4,4,3,3,4,5,3,3,5,3,2,2,3,5,4,2,1,5,3,3,3,2,2,1,5,3,5,5,2,5,2,4,4,1,1,5,3,1,4,5,2,3,5,5,1,2,1,3,3,2,2,2,1,1,2,2,3,3,3,2,3,2)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
